# The document contains three places where an "<id>...</id>" tag was
# originally split across three separately-formatted runs, e.g.:
#   run1 (Courier New): "<id>"
#   run2 (Arial):        "p170v_1"
#   run3 (Courier New):  "</id>"
#
# The edit collapses each of these three runs into a single run
# (using the formatting of the first/opening run) containing the full
# "<id>p170v_N</id>" text, for N = 1, 2, 3.
#
# A Find/Replace whose search text spans the run boundaries causes Word
# to rewrite the matched range as a single run using the formatting of
# the first run in the match - exactly the behavior we need.

$d = $word.ActiveDocument

$ids = @("p170v_1", "p170v_2", "p170v_3")

foreach ($pid in $ids) {
    $target = "<id>" + $pid + "</id>"

    $rng = $d.Content
    $found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, `
                                $true, 1, $false, $target, 2)

    if (-not $found) {
        Write-Host "WARNING: could not find '$target'"
    }
}
